$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.938.02"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.370.04"
$ws.Range("E3").Value = "  +4.00%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'235.83"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'0.661"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").Value = "'73.22"
$ws.Range("E7").Value = "  +13.64%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  +21.58%  "
$ws.Range("D10").Value = "'0.0991"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("D11").Value = "'28.27"
$ws.Range("E11").Value = "  +6.86%  "
$ws.Range("D12").Value = "2.727.77"
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("D14").Value = "'16.98"
$ws.Range("E14").Value = "  +12.45%  "
$ws.Range("D15").Value = "'6.69"
$ws.Range("E15").Value = "  +9.95%  "
$ws.Range("D16").Value = "'0.886"
$ws.Range("E16").Value = "  +7.43%  "
$ws.Range("D17").Value = "2.370.91"
$ws.Range("E17").Value = "  +4.29%  "
$ws.Range("D18").Value = "43.900.31"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("D20").Value = "'76.26"
$ws.Range("E20").Value = "  +4.45%  "
$ws.Range("D21").Value = "'6.37"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").Value = "'251.65"
$ws.Range("E22").Value = "  +1.56%  "
$ws.Range("B23").Value = "WEMIXToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D23").Value = "'3.80"
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'2.50"
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").Value = "'10.36"
$ws.Range("E26").Value = "  +6.63%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").Value = "'22.56"
$ws.Range("E28").Value = "  +3.55%  "
$ws.Range("D29").Value = "'173.43"
$ws.Range("E29").Value = "  -0.19%  "
$ws.Range("D30").Value = "'1.55"
$ws.Range("E30").Value = "  +8.90%  "
$ws.Range("D31").Value = "'0.133"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("D33").Value = "'5.22"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("D34").Value = "'0.0710"
$ws.Range("E34").Value = "  +4.67%  "
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("D36").Value = "'3.77"
$ws.Range("E36").Value = "  +5.00%  "
$ws.Range("D37").Value = "'2.45"
$ws.Range("E37").Value = "  +6.84%  "
$ws.Range("D38").Value = "'6.44"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "'0.0266"
$ws.Range("E39").Value = "  +6.45%  "
$ws.Range("D40").Value = "'19.55"
$ws.Range("E40").Value = "  +13.76%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.96"
$ws.Range("E41").Value = "  +2.42%  "
$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("E43").Value = "  +8.91%  "
$ws.Range("D44").Value = "'1.22"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.0969"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'98.53"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  +13.92%  "
$ws.Range("D49").Value = "1.443.12"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'2.32"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.597.54"
$ws.Range("E51").Value = "  +4.34%  "
